# Update the "repaymentstrategy" input on the ProductLoanInput sheet from
# "RBI (India)" to "Overdue/Due Fee/Int,Principal" and give it its own
# distinct (Arial 10 / green fill) formatting, matching the style already
# used for the other green input cells on the sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")

$cell = $ws1.Range("B17")
$cell.Value = "Overdue/Due Fee/Int,Principal"
$cell.Font.Name = "Arial"
$cell.Font.Size = 10
$cell.Interior.Color = 5296274

# Move the active selection to the cell that was just edited.
$ws1.Range("B17").Select()
